$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "27.28") are stored as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.707.38'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '2.475.28'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '576.46'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('D6').Value = '149.14'
$ws.Range('E6').Value = '  +2.63%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +1.84%  '
$ws.Range('D9').Value = '2.471.95'
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('D14').Value = '27.28'
$ws.Range('E14').Value = '  +1.72%  '
$ws.Range('D15').Value = '0.0000182'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '2.921.31'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '63.486.80'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').Value = '2.482.24'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('D19').Value = '11.49'
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('D20').Value = '7.44'
$ws.Range('E20').Value = '  +7.47%  '
$ws.Range('D21').Value = '331.09'
$ws.Range('E21').Value = '  +2.38%  '
$ws.Range('D22').Value = '4.23'
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('D23').Value = '2.13'
$ws.Range('E23').Value = '  +18.65%  '
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '66.04'
$ws.Range('E25').Value = '  -1.86%  '
$ws.Range('D26').Value = '630.62'
$ws.Range('E26').Value = '  +11.31%  '
$ws.Range('D27').Value = '9.22'
$ws.Range('E27').Value = '  +5.96%  '
$ws.Range('E28').Value = '  +3.91%  '
$ws.Range('E29').Value = '  +6.27%  '
$ws.Range('D30').Value = '2.603.35'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').Value = '8.43'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('E33').Value = '  -2.27%  '
$ws.Range('D34').Value = '1.92'
$ws.Range('E34').Value = '  +1.62%  '
$ws.Range('D35').Value = '5.28'
$ws.Range('E35').Value = '  +8.49%  '
$ws.Range('D36').Value = '1.56'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '0.384'
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('D39').Value = '5.52'
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D40').Value = '18.90'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').Value = '2.76'
$ws.Range('E41').Value = '  +14.35%  '
$ws.Range('D42').Value = '147.33'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('D45').Value = '151.13'
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').Value = '3.79'
$ws.Range('E46').Value = '  +3.20%  '
$ws.Range('D47').Value = '21.60'
$ws.Range('E47').Value = '  +5.43%  '
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('D49').Value = '0.608'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').Value = '0.0237'
$ws.Range('E50').Value = '  +2.38%  '
$ws.Range('D51').Value = '0.0921'
$ws.Range('E51').Value = '  -0.58%  '
